# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel auto-converting
# numeric-looking strings like "484.67" or "1.00" into real numbers), then reset the
# cell style back to Normal so no stray NumberFormat/quote-prefix style lingers on it.
function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Plain text values (names/URLs/already non-numeric-shaped strings) can be set directly.
function Set-PlainValue {
    param($Address, $Text)
    $ws.Range($Address).Value = $Text
}

Set-PlainValue "D2" "68.343.16"
Set-PlainValue "E2" "  +1.57%  "
Set-PlainValue "D3" "3.925.92"
Set-PlainValue "E3" "  -0.55%  "
Set-PlainValue "E4" "  +0.02%  "
Set-TextValue "D5" "484.67"
Set-PlainValue "E5" "  +2.86%  "
Set-TextValue "D6" "147.25"
Set-PlainValue "E6" "  +0.97%  "
Set-TextValue "D7" "0.626"
Set-PlainValue "E7" "  +0.33%  "
Set-PlainValue "E8" "  -0.01%  "
Set-TextValue "D9" "0.732"
Set-PlainValue "E9" "  -0.17%  "
Set-TextValue "D10" "0.167"
Set-PlainValue "E10" "  +1.67%  "
Set-TextValue "D11" "0.0000351"
Set-PlainValue "E11" "  +4.01%  "
Set-TextValue "D12" "42.83"
Set-PlainValue "E12" "  -1.51%  "
Set-TextValue "D13" "10.70"
Set-PlainValue "E13" "  +2.89%  "
Set-PlainValue "D14" "4.557.31"
Set-PlainValue "E14" "  -0.39%  "
Set-PlainValue "B15" "WrappedEther"
Set-PlainValue "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-PlainValue "D15" "3.955.19"
Set-PlainValue "E15" "  +1.54%  "
Set-PlainValue "B16" "Uniswap"
Set-PlainValue "C16" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D16" "14.43"
Set-PlainValue "E16" "  -5.09%  "
Set-PlainValue "E17" "  -0.65%  "
Set-TextValue "D18" "19.92"
Set-PlainValue "E18" "  +0.31%  "
Set-TextValue "D19" "1.13"
Set-PlainValue "E19" "  -2.20%  "
Set-PlainValue "D20" "68.448.94"
Set-PlainValue "E20" "  +1.32%  "
Set-TextValue "D21" "441.28"
Set-PlainValue "E21" "  +1.86%  "
Set-TextValue "D22" "3.49"
Set-PlainValue "E22" "  +3.16%  "
Set-TextValue "D23" "15.04"
Set-PlainValue "E23" "  +3.38%  "
Set-TextValue "D24" "88.37"
Set-PlainValue "E24" "  +1.11%  "
Set-PlainValue "B25" "Filecoin"
Set-PlainValue "C25" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D25" "11.53"
Set-PlainValue "E25" "  +11.82%  "
Set-PlainValue "B26" "RenderToken"
Set-PlainValue "C26" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D26" "11.20"
Set-PlainValue "E26" "  +16.00%  "
Set-TextValue "D27" "3.62"
Set-PlainValue "E27" "  +0.83%  "
Set-TextValue "D28" "39.01"
Set-PlainValue "E28" "  +1.12%  "
Set-TextValue "D29" "5.81"
Set-PlainValue "E29" "  +0.93%  "
Set-TextValue "D30" "717.09"
Set-PlainValue "E30" "  -1.10%  "
Set-TextValue "D31" "13.72"
Set-PlainValue "E31" "  +0.78%  "
Set-PlainValue "E32" "  -1.88%  "
Set-PlainValue "E33" "  +3.74%  "
Set-PlainValue "B34" "NEARProtocol"
Set-PlainValue "C34" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D34" "6.24"
Set-PlainValue "E34" "  +15.96%  "
Set-PlainValue "B35" "InjectiveProtocol"
Set-PlainValue "C35" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D35" "42.07"
Set-PlainValue "E35" "  -2.55%  "
Set-PlainValue "E36" "  +13.42%  "
Set-TextValue "D37" "60.96"
Set-PlainValue "E37" "  +5.25%  "
Set-TextValue "D38" "0.149"
Set-PlainValue "E38" "  -1.93%  "
Set-TextValue "D39" "0.399"
Set-PlainValue "E39" "  +18.37%  "
Set-TextValue "D40" "1.00"
Set-PlainValue "E40" "  +0.06%  "
Set-TextValue "D41" "2.95"
Set-PlainValue "E41" "  +14.74%  "
Set-TextValue "D42" "3.22"
Set-PlainValue "E42" "  +5.80%  "
Set-TextValue "D43" "0.0481"
Set-PlainValue "E43" "  +0.65%  "
Set-TextValue "D44" "2.94"
Set-PlainValue "E44" "  +4.28%  "
Set-PlainValue "E45" "  +0.61%  "
Set-PlainValue "E46" "  +0.11%  "
Set-PlainValue "D47" "0.0₆0357"
Set-PlainValue "E47" "  +36.62%  "
Set-TextValue "D48" "3.27"
Set-PlainValue "E48" "  +3.43%  "
Set-PlainValue "E49" "  -1.50%  "
Set-TextValue "D50" "2.14"
Set-PlainValue "E50" "  -2.16%  "
Set-TextValue "D51" "145.81"
Set-PlainValue "E51" "  -0.90%  "
